$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 145
$ws.Range("I12").Value = 90
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 90
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = -540

$ws.Range("H96").Value = 1458.3334
$ws.Range("I96").Value = 1362.5
$ws.Range("J96").Value = 1650
$ws.Range("K96").Value = 4087.5
$ws.Range("L96").Value = 4950
$ws.Range("M96").Value = -2714.5
$ws.Range("N96").Value = -7696

$ws.Range("H125").Value = 1942.9333
$ws.Range("I125").Value = 854.6667
$ws.Range("J125").Value = 2215
$ws.Range("K125").Value = 7692.0003
$ws.Range("L125").Value = 19935
$ws.Range("M125").Value = -5232.0003
$ws.Range("N125").Value = -24855

$ws.Range("H138").Value = 2340.12
$ws.Range("J138").Value = 2415.5532
$ws.Range("L138").Value = 7246.659599999999
$ws.Range("N138").Value = -17526.6596

$ws.Range("H140").Value = 53378.8
$ws.Range("J140").Value = 53378.8
$ws.Range("L140").Value = 53378.8
$ws.Range("N140").Value = -63738.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1075959.5
$ws.Range("I32").Value = 1140344.1
$ws.Range("K32").Value = 1140344.1
$ws.Range("M32").Value = -1140057.1

$ws.Range("H132").Value = 1511560.9
$ws.Range("I132").Value = 2163.7144
$ws.Range("J132").Value = 4813367
$ws.Range("K132").Value = 6491.1432
$ws.Range("L132").Value = 14440101
$ws.Range("M132").Value = -3961.1432
$ws.Range("N132").Value = -14445161

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3459.5715
$ws.Range("I80").Value = 4383.4
$ws.Range("K80").Value = 4383.4
$ws.Range("M80").Value = -3385.4

$ws.Range("H83").Value = 3459.5715
$ws.Range("I83").Value = 4383.4
$ws.Range("K83").Value = 21917
$ws.Range("M83").Value = -16925

$ws.Range("H134").Value = 3316.9355
$ws.Range("I134").Value = 3188.5217
$ws.Range("J134").Value = 3686.125
$ws.Range("K134").Value = 9565.5651
$ws.Range("L134").Value = 11058.375
$ws.Range("M134").Value = -7030.5651
$ws.Range("N134").Value = -16128.375

$ws.Range("H135").Value = 39900
$ws.Range("J135").Value = 39900
$ws.Range("L135").Value = 39900
$ws.Range("N135").Value = -50040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7148.622
$ws.Range("I31").Value = 2194.6316
$ws.Range("K31").Value = 2194.6316
$ws.Range("M31").Value = -1899.6316

$ws.Range("H34").Value = 7148.622
$ws.Range("I34").Value = 2194.6316
$ws.Range("K34").Value = 2194.6316
$ws.Range("M34").Value = -1992.6316

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 396.875
$ws.Range("I26").Value = 75.375
$ws.Range("J26").Value = 504.04166
$ws.Range("K26").Value = 226.125
$ws.Range("L26").Value = 1512.12498
$ws.Range("M26").Value = 61.875
$ws.Range("N26").Value = -2088.12498

$ws.Range("H122").Value = 3208.1
$ws.Range("I122").Value = 462.75
$ws.Range("J122").Value = 5038.3335
$ws.Range("K122").Value = 4164.75
$ws.Range("L122").Value = 45345.0015
$ws.Range("M122").Value = -1714.75
$ws.Range("N122").Value = -50245.0015

$ws.Range("H131").Value = 3793.9048
$ws.Range("I131").Value = 293.16666
$ws.Range("J131").Value = 4377.3613
$ws.Range("K131").Value = 879.4999799999999
$ws.Range("L131").Value = 13132.0839
$ws.Range("M131").Value = 4160.50002
$ws.Range("N131").Value = -23212.0839

$ws.Range("H137").Value = 101389.27
$ws.Range("I137").Value = 13618.625
$ws.Range("K137").Value = 40855.875
$ws.Range("M137").Value = -35755.875

$ws.Range("H139").Value = 4192.564
$ws.Range("J139").Value = 4850.276
$ws.Range("L139").Value = 14550.828
$ws.Range("N139").Value = -24830.828

$ws.Range("H140").Value = 1946.8667
$ws.Range("I140").Value = 1677
$ws.Range("K140").Value = 5031
$ws.Range("M140").Value = 149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 481
$ws.Range("I107").Value = 526.8
$ws.Range("J107").Value = 366.5
$ws.Range("K107").Value = 526.8
$ws.Range("L107").Value = 366.5
$ws.Range("M107").Value = 1393.2
$ws.Range("N107").Value = -4206.5

$ws.Range("H132").Value = 3080.389
$ws.Range("I132").Value = 2152.6365
$ws.Range("J132").Value = 4538.2856
$ws.Range("K132").Value = 6457.9095
$ws.Range("L132").Value = 13614.8568
$ws.Range("M132").Value = -3927.9095
$ws.Range("N132").Value = -18674.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3195.1428
$ws.Range("I16").Value = 3217.25
$ws.Range("J16").Value = 3165.6667
$ws.Range("K16").Value = 3217.25
$ws.Range("L16").Value = 3165.6667
$ws.Range("M16").Value = -3047.25
$ws.Range("N16").Value = -3505.6667

$ws.Range("H46").Value = 1561.5
$ws.Range("I46").Value = 1123
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1123
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -935
$ws.Range("N46").Value = -2376

$ws.Range("H55").Value = 348.3
$ws.Range("I55").Value = 255.5
$ws.Range("J55").Value = 564.8333
$ws.Range("K55").Value = 255.5
$ws.Range("L55").Value = 564.8333
$ws.Range("M55").Value = -82.5
$ws.Range("N55").Value = -910.8333

$ws.Range("H68").Value = 1897
$ws.Range("I68").Value = 1994
$ws.Range("J68").Value = 1800
$ws.Range("K68").Value = 1994
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -1245
$ws.Range("N68").Value = -3298

$ws.Range("H71").Value = 1897
$ws.Range("I71").Value = 1994
$ws.Range("J71").Value = 1800
$ws.Range("K71").Value = 9970
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -6226
$ws.Range("N71").Value = -16488

$ws.Range("H132").Value = 3473.027
$ws.Range("I132").Value = 3163.1853
$ws.Range("J132").Value = 4309.6
$ws.Range("K132").Value = 9489.555899999999
$ws.Range("L132").Value = 12928.8
$ws.Range("M132").Value = -6959.555899999999
$ws.Range("N132").Value = -17988.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1085.1305
$ws.Range("I107").Value = 1068.8334
$ws.Range("K107").Value = 3206.5002
$ws.Range("M107").Value = -1286.5002

$ws.Range("H132").Value = 3475287.5
$ws.Range("I132").Value = 4028.5
$ws.Range("J132").Value = 5684270.5
$ws.Range("K132").Value = 12085.5
$ws.Range("L132").Value = 17052811.5
$ws.Range("M132").Value = -9555.5
$ws.Range("N132").Value = -17057871.5
